$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D2").Value = 42776

$ws.Range("D3").Value = "+"
$ws.Range("D4").Value = "+"
$ws.Range("D5").Value = "+"
$ws.Range("D6").Value = "+"
$ws.Range("D7").Value = "-"
$ws.Range("D8").Value = "-"
$ws.Range("D9").Value = "+"
$ws.Range("D10").Value = "-"
$ws.Range("D11").Value = "+"
$ws.Range("D12").Value = "+"
$ws.Range("D13").Value = "+"
$ws.Range("D14").Value = "+"
$ws.Range("D15").Value = "+"
$ws.Range("D16").Value = "+"
$ws.Range("D17").Value = "+"
$ws.Range("D18").Value = "+"
$ws.Range("D19").Value = "+"
$ws.Range("D20").Value = "-"
$ws.Range("D21").Value = "+"
$ws.Range("D22").Value = "+"
$ws.Range("D23").Value = "-"

$ws.Columns.Item(4).ColumnWidth = 9.25

$ws.Range("N12").Select()
